{"js": "const NEW_DATE = \"2024-12-11 Wednesday\";\nconst NEW_CELLS = [\n  [\"74-9=\", \"40+30=\", \"2+72=\", \"83-11=\", \"6+80=\"],\n  [\"46+18=\", \"20+53=\", \"89-78=\", \"34-8=\", \"10+40=\"],\n  [\"33+58=\", \"64+17=\", \"15-14=\", \"86-4=\", \"46+46=\"],\n  [\"53-4=\", \"55-12=\", \"79-58=\", \"22-18=\", \"36+23=\"],\n  [\"3+3=\", \"27+15=\", \"79-72=\", \"57+12=\", \"51+31=\"],\n  [\"37-37=\", \"24+30=\", \"61-40=\", \"66-52=\", \"1+45=\"],\n  [\"42-10=\", \"68-35=\", \"11+65=\", \"69+10=\", \"99-44=\"],\n  [\"11+75=\", \"71-2=\", \"73-34=\", \"16+33=\", \"13+67=\"],\n  [\"48+3=\", \"87-86=\", \"65-9=\", \"87-41=\", \"95-5=\"],\n  [\"9+10=\", \"29+68=\", \"95-57=\", \"40+29=\", \"61+36=\"],\n  [\"81-57=\", \"81-8=\", \"42+10=\", \"63-35=\", \"26+22=\"],\n  [\"84-53=\", \"17+73=\", \"57-43=\", \"10+12=\", \"39+18=\"],\n  [\"93-55=\", \"34-8=\", \"80-5=\", \"78-63=\", \"93-74=\"],\n  [\"7+49=\", \"61+4=\", \"73-3=\", \"10+8=\", \"41+34=\"],\n  [\"0+36=\", \"33+61=\", \"60+11=\", \"6+62=\", \"3+32=\"],\n  [\"40-17=\", \"99-76=\", \"29+40=\", \"99-12=\", \"61-21=\"],\n  [\"73-34=\", \"85-11=\", \"98-4=\", \"16+55=\", \"10+51=\"],\n  [\"77-70=\", \"26-9=\", \"29-4=\", \"24+7=\", \"36-15=\"],\n  [\"27-15=\", \"57+39=\", \"60-5=\", \"15+11=\", \"74+13=\"],\n  [\"66+25=\", \"75-36=\", \"48-20=\", \"72+0=\", \"91-18=\"],\n];\n\n\n// 1) Update the header date paragraph (first paragraph in the body, outside the table).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.text.includes(\"2024-12-10 Tuesday\")) {\n    p.insertText(NEW_DATE, \"Replace\");\n    break;\n  }\n}\n\n// 2) Update every arithmetic-expression cell in the (single) table, row by row,\n//    left to right, matching the document's row/column order.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let r = 0; r < rows.items.length; r++) {\n  const cells = rows.items[r].cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (let c = 0; c < cells.items.length; c++) {\n    cells.items[c].value = NEW_CELLS[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "$NewDate = \"2024-12-11 Wednesday\"\n$NewCells = @(\n    @(\"74-9=\", \"40+30=\", \"2+72=\", \"83-11=\", \"6+80=\"),\n    @(\"46+18=\", \"20+53=\", \"89-78=\", \"34-8=\", \"10+40=\"),\n    @(\"33+58=\", \"64+17=\", \"15-14=\", \"86-4=\", \"46+46=\"),\n    @(\"53-4=\", \"55-12=\", \"79-58=\", \"22-18=\", \"36+23=\"),\n    @(\"3+3=\", \"27+15=\", \"79-72=\", \"57+12=\", \"51+31=\"),\n    @(\"37-37=\", \"24+30=\", \"61-40=\", \"66-52=\", \"1+45=\"),\n    @(\"42-10=\", \"68-35=\", \"11+65=\", \"69+10=\", \"99-44=\"),\n    @(\"11+75=\", \"71-2=\", \"73-34=\", \"16+33=\", \"13+67=\"),\n    @(\"48+3=\", \"87-86=\", \"65-9=\", \"87-41=\", \"95-5=\"),\n    @(\"9+10=\", \"29+68=\", \"95-57=\", \"40+29=\", \"61+36=\"),\n    @(\"81-57=\", \"81-8=\", \"42+10=\", \"63-35=\", \"26+22=\"),\n    @(\"84-53=\", \"17+73=\", \"57-43=\", \"10+12=\", \"39+18=\"),\n    @(\"93-55=\", \"34-8=\", \"80-5=\", \"78-63=\", \"93-74=\"),\n    @(\"7+49=\", \"61+4=\", \"73-3=\", \"10+8=\", \"41+34=\"),\n    @(\"0+36=\", \"33+61=\", \"60+11=\", \"6+62=\", \"3+32=\"),\n    @(\"40-17=\", \"99-76=\", \"29+40=\", \"99-12=\", \"61-21=\"),\n    @(\"73-34=\", \"85-11=\", \"98-4=\", \"16+55=\", \"10+51=\"),\n    @(\"77-70=\", \"26-9=\", \"29-4=\", \"24+7=\", \"36-15=\"),\n    @(\"27-15=\", \"57+39=\", \"60-5=\", \"15+11=\", \"74+13=\"),\n    @(\"66+25=\", \"75-36=\", \"48-20=\", \"72+0=\", \"91-18=\"),\n)\n\n\n$d = $word.ActiveDocument\n\n# 1) Update the header date paragraph (first paragraph, before the table).\n$d.Paragraphs.Item(1).Range.Text = $NewDate\n\n# 2) Update every arithmetic-expression cell in the (single) table, row by\n#    row, left to right, matching the document's row/column order.\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $NewCells[$r - 1][$c - 1]\n    }\n}\n"}
